$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the expiry date value in C2 from "03/20" to "09/20"
$ws.Range("C2").Value = "09/20"

# Update the selected cell/range to reflect the new selection (C6)
$ws.Range("C6").Select()
